# Weekly update: insert a new data row at the top of the Piña (Vega Modelo de
# Temuco) table, pushing all existing rows down by one.
#
# Before the edit, data occupied rows 2-819 (dimension A1:T819). The new
# record becomes the new row 733 and every row that used to be 733-819
# shifts down to 734-820 (dimension becomes A1:T820).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 733, shifting rows 733:819
# down to 734:820.
$ws.Rows.Item(733).Insert()

# Populate the newly inserted row 733 with the new weekly data point.
$ws.Cells.Item(733, 1).Value = 10
$ws.Cells.Item(733, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(733, 3).Value = "La Araucanía"
$ws.Cells.Item(733, 4).Value = 45212
$ws.Cells.Item(733, 5).Value = 9
$ws.Cells.Item(733, 6).Value = "Fruta"
$ws.Cells.Item(733, 7).Value = 100108
$ws.Cells.Item(733, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(733, 9).Value = 100108005
$ws.Cells.Item(733, 10).Value = "Piña"
$ws.Cells.Item(733, 11).Value = "Caramelo"
$ws.Cells.Item(733, 12).Value = "Primera"
$ws.Cells.Item(733, 13).Value = 100
$ws.Cells.Item(733, 14).Value = 25000
$ws.Cells.Item(733, 15).Value = 25000
$ws.Cells.Item(733, 16).Value = 25000
$ws.Cells.Item(733, 17).Value = "`$/caja 12 unidades"
$ws.Cells.Item(733, 18).Value = "Ecuador"
$ws.Cells.Item(733, 19).Value = 2083
$ws.Cells.Item(733, 20).Value = 12
